# Applies the cryptos price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'20.338.53"
$ws.Range("E2").Value = "  +1.85%  "

$ws.Range("D3").Value = "'1.454.28"
$ws.Range("E3").Value = "  +3.13%  "

$ws.Range("D4").Value = "'1.015"
$ws.Range("E4").Value = "  +1.37%  "

$ws.Range("D5").Value = "'278.47"
$ws.Range("E5").Value = "  +1.90%  "

$ws.Range("D6").Value = "'0.8961"
$ws.Range("E6").Value = "  -10.43%  "

$ws.Range("D7").Value = "'0.3689"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'0.3147"
$ws.Range("E8").Value = "  +2.62%  "

$ws.Range("D9").Value = "'39.17"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").Value = "'1.026"
$ws.Range("E10").Value = "  +2.78%  "

$ws.Range("D11").Value = "'0.06514"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "'1.007"
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13").Value = "'5.433"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").Value = "'17.54"
$ws.Range("E14").Value = "  +3.13%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.119"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'1.461.71"
$ws.Range("E16").Value = "  +3.82%  "

$ws.Range("D17").Value = "'0.00001019"
$ws.Range("E17").Value = "  +1.36%  "

$ws.Range("D18").Value = "'0.05610"
$ws.Range("E18").Value = "  -2.64%  "

$ws.Range("D19").Value = "'0.9032"
$ws.Range("E19").Value = "  -9.67%  "

$ws.Range("D20").Value = "'67.75"
$ws.Range("E20").Value = "  -7.83%  "

$ws.Range("D21").Value = "'5.474"
$ws.Range("E21").Value = "  -2.26%  "

$ws.Range("D22").Value = "'14.44"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'11.10"
$ws.Range("E23").Value = "  +2.46%  "

$ws.Range("D24").Value = "'2.259"
$ws.Range("E24").Value = "  -2.34%  "

$ws.Range("D25").Value = "'20.444.62"
$ws.Range("E25").Value = "  +2.34%  "

$ws.Range("D26").Value = "'2.207"
$ws.Range("E26").Value = "  -2.50%  "

$ws.Range("D27").Value = "'135.44"
$ws.Range("E27").Value = "  -2.22%  "

$ws.Range("D28").Value = "'17.08"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Range("D29").Value = "'1.622.31"
$ws.Range("E29").Value = "  +3.55%  "

$ws.Range("D30").Value = "'111.09"
$ws.Range("E30").Value = "  +2.01%  "

$ws.Range("D31").Value = "'3.635"
$ws.Range("E31").Value = "  -4.97%  "

$ws.Range("D32").Value = "'0.8115"
$ws.Range("E32").Value = "  -4.40%  "

$ws.Range("D33").Value = "'4.918"
$ws.Range("E33").Value = "  -8.45%  "

$ws.Range("D34").Value = "'0.07694"

$ws.Range("D35").Value = "'0.05987"
$ws.Range("E35").Value = "  +3.28%  "

$ws.Range("D36").Value = "'1.420"
$ws.Range("E36").Value = "  +10.92%  "

$ws.Range("D37").Value = "'4.742"
$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("E38").Value = "  +7.59%  "

$ws.Range("D39").Value = "'0.02026"
$ws.Range("E39").Value = "  -0.98%  "

$ws.Range("D40").Value = "'10.31"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").Value = "'0.1840"
$ws.Range("E41").Value = "  -4.57%  "

$ws.Range("D42").Value = "'0.9162"
$ws.Range("E42").Value = "  -8.37%  "

$ws.Range("D43").Value = "'3.562"
$ws.Range("E43").Value = "  +0.87%  "

$ws.Range("D44").Value = "'0.5280"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'6.801"
$ws.Range("E45").Value = "  -19.52%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'11.91"
$ws.Range("E46").Value = "  -1.97%  "

$ws.Range("D47").Value = "'121.32"
$ws.Range("E47").Value = "  +10.30%  "

$ws.Range("D48").Value = "'0.5160"
$ws.Range("E48").Value = "  +0.95%  "

$ws.Range("D49").Value = "'1.773"
$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("D50").Value = "'0.06361"
$ws.Range("E50").Value = "  +3.16%  "

$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  +0.04%  "

